$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.227.19"
$ws.Range("E2").Value = "  -5.94%  "

$ws.Range("D3").Value = "2.222.59"
$ws.Range("E3").Value = "  -5.96%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.16"
$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("E6").Value = "  -5.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.21"
$ws.Range("E7").Value = "  -6.07%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.557"
$ws.Range("E9").Value = "  -7.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.21"
$ws.Range("E10").Value = "  +5.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  -7.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.30"
$ws.Range("E12").Value = "  -2.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  -3.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.76"
$ws.Range("E14").Value = "  -7.74%  "

$ws.Range("D15").Value = "2.552.56"
$ws.Range("E15").Value = "  -5.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.83"
$ws.Range("E16").Value = "  -9.91%  "

$ws.Range("E17").Value = "  -8.86%  "

$ws.Range("D18").Value = "2.222.12"
$ws.Range("E18").Value = "  -6.25%  "

$ws.Range("D19").Value = "41.309.84"
$ws.Range("E19").Value = "  -5.66%  "

$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -8.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.39"
$ws.Range("E21").Value = "  -6.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  -8.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.93"
$ws.Range("E23").Value = "  -8.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  +12.64%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.67"
$ws.Range("E26").Value = "  -3.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.43"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.86"
$ws.Range("E28").Value = "  -7.56%  "

$ws.Range("E29").Value = "  -4.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.40"
$ws.Range("E30").Value = "  -2.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.59"
$ws.Range("E31").Value = "  -8.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -8.14%  "

$ws.Range("E33").Value = "  -7.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0715"
$ws.Range("E34").Value = "  -5.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.23"
$ws.Range("E35").Value = "  -4.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  -10.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.90"
$ws.Range("E37").Value = "  +2.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.00"
$ws.Range("E38").Value = "  +14.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0276"
$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("E40").Value = "  -5.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.86"
$ws.Range("E41").Value = "  -11.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.61"
$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("E43").Value = "  +2.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.00"
$ws.Range("E44").Value = "  -11.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.95"
$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.100"
$ws.Range("E46").Value = "  -6.87%  "

$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.64"
$ws.Range("E47").Value = "  +6.28%  "

$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.78"
$ws.Range("E48").Value = "  +10.40%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("E50").Value = "  -6.10%  "

$ws.Range("B51").Value = "BitTorrent-New"
$ws.Range("C51").Value = "https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt"
$ws.Range("D51").Value = "0.0₃0147"
$ws.Range("E51").Value = "  +9.74%  "
